$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finish out row 10 (Holiday week 10/15-19) - A10/B10/C10 already filled in
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

# Row 11: week of 10/22-26
$ws.Range("A11").Value = "10/22-26"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 2.5
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 4

# Row 12: week of 10/29-11/2 (no Friday hours logged)
$ws.Range("A12").Value = "10/29-11/2"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 6

# Extend the "Total" column formula down through the new rows
$ws.Range("G2:G12").Formula = "=SUM(B2:F2)"

$ws.Range("G12").Select()
